# Fix the mislabeled "2050" (or "2041-2050") column header that had been
# accidentally overwritten with a stray numeric value (712.444928568978),
# and remove the "Total" row that used to sit at the bottom of several
# tables.

$wb = $excel.ActiveWorkbook

function Set-HeaderLabel($ws, $cellAddress, $text) {
    # Write the text via a literal-string formula so Excel treats it as
    # text even though it looks numeric (e.g. "2050"), then convert the
    # formula to a plain value in place. This avoids turning the label
    # into a real number while leaving the cell's existing formatting
    # (bold, centered, bordered header style) untouched.
    $cell = $ws.Range($cellAddress)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $ws.Application.CutCopyMode = 0
}

# Sheets whose column E header (row 1) must read "2050" and which have a
# trailing "Total" row (row 13) that must be removed.
$simpleSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)"
)

foreach ($name in $simpleSheets) {
    $ws = $wb.Worksheets.Item($name)

    Set-HeaderLabel $ws "E1" "2050"

    # Drop the "Total" row.
    $ws.Rows.Item(13).Delete()
}

# "Potencia Incremental - SIN(MW)" follows the same pattern, but its
# column headers are period ranges, so E1 becomes "2041-2050".
$wsInc = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-HeaderLabel $wsInc "E1" "2041-2050"
$wsInc.Rows.Item(13).Delete()

# "Emissoes Totais (MtCO2eq)" only needs the mislabeled header fixed; it
# never had a "Total" row.
$wsEm = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-HeaderLabel $wsEm "E1" "2050"

# "Custo Total (bilhões de R$)" keeps its header as-is, but also loses its
# trailing "Total" row (row 4).
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Rows.Item(4).Delete()
